# Add "Mid Paper 1" / "Mid Paper 2" columns (F, G) to both the
# "Senior Six" and "Senior Five" sheets, mirroring the existing
# "Paper 1" / "Paper 2" (D, E) marks into the new columns, and leave
# "Senior Six" as the active / selected sheet with F1:G6 selected
# (matching the selection that was also left on "Senior Five").

$wb = $excel.ActiveWorkbook

$sheetNames = @("Senior Six", "Senior Five")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # New header cells.
    $ws.Range("F1").Value2 = "Mid Paper 1"
    $ws.Range("G1").Value2 = "Mid Paper 2"

    # Mirror the Paper 1 / Paper 2 marks already present for the
    # students that have marks recorded (rows 2, 5 and 6).
    $ws.Range("F2").Value2 = $ws.Range("D2").Value2
    $ws.Range("G2").Value2 = $ws.Range("E2").Value2

    $ws.Range("F5").Value2 = $ws.Range("D5").Value2
    $ws.Range("G5").Value2 = $ws.Range("E5").Value2

    $ws.Range("F6").Value2 = $ws.Range("D6").Value2
    $ws.Range("G6").Value2 = $ws.Range("E6").Value2
}

# Leave the same F1:G6 selection on "Senior Five" ...
$wsFive = $wb.Worksheets.Item("Senior Five")
$wsFive.Activate() | Out-Null
$wsFive.Range("F1:G6").Select() | Out-Null

# ... but make "Senior Six" the active / selected tab, also with
# F1:G6 selected.
$wsSix = $wb.Worksheets.Item("Senior Six")
$wsSix.Activate() | Out-Null
$wsSix.Range("F1:G6").Select() | Out-Null
